$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Untagged")
Write-Host ($ws.Columns.Item(13).ColumnWidth)
